$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = '0.0008939121986040734'
$ws.Range("C3").Value = '0.0006749952681096833'
$ws.Range("C4").Value = '0.000417295443565566'
$ws.Range("C5").Value = '0.0003780690076164485'
$ws.Range("C6").Value = '0.0003190170992704266'
$ws.Range("C7").Value = '0.000290456345059997'
$ws.Range("C8").Value = '0.0002642698241043285'
$ws.Range("C9").Value = '0.000178032412925192'
$ws.Range("C10").Value = '0.0001610359535311469'
$ws.Range("C11").Value = '0.0001198143315554234'
$ws.Range("C12").Value = '0.0001141978222358567'
$ws.Range("C13").Value = '0.0001029718203552676'
$ws.Range("C14").Value = '8.983421693672707e-05'
$ws.Range("C15").Value = '8.477129753548767e-05'
$ws.Range("C16").Value = '7.814359288473629e-05'
$ws.Range("C17").Value = '6.833080695786403e-05'
$ws.Range("C18").Value = '6.689694169855833e-05'
$ws.Range("C19").Value = '6.153987659244439e-05'
$ws.Range("C20").Value = '6.09994325209579e-05'
$ws.Range("C21").Value = '6.031290888366295e-05'
$ws.Range("C22").Value = '4.939792260511538e-05'
$ws.Range("C23").Value = '4.65261629682723e-05'
$ws.Range("A24").Value = 11
$ws.Range("B24").Value = '(h$_{pr}$)$_{3}$'
$ws.Range("C24").Value = '4.502228883294549e-05'
$ws.Range("A25").Value = 6
$ws.Range("B25").Value = '(h$_{pq}$)$_{2}$'
$ws.Range("C25").Value = '4.467880786697934e-05'
$ws.Range("C26").Value = '4.342032329911382e-05'
$ws.Range("C27").Value = '4.000277612047268e-05'
$ws.Range("A28").Value = 42
$ws.Range("B28").Value = '$F_{s}$'
$ws.Range("C28").Value = '3.959587678930198e-05'
$ws.Range("A29").Value = 16
$ws.Range("B29").Value = '(h$_{r}$)$_{2}$'
$ws.Range("C29").Value = '3.834311846226618e-05'
$ws.Range("C30").Value = '3.794023953174323e-05'
$ws.Range("C31").Value = '3.792331513884148e-05'
$ws.Range("C32").Value = '3.775333273441231e-05'
$ws.Range("C33").Value = '3.539565223401197e-05'
$ws.Range("C34").Value = '3.193173444407881e-05'
$ws.Range("C35").Value = '2.962519600775321e-05'
$ws.Range("C36").Value = '2.773982478672589e-05'
$ws.Range("C37").Value = '2.758421325778318e-05'
$ws.Range("C38").Value = '2.742261031145537e-05'
$ws.Range("C39").Value = '2.688817442583975e-05'
$ws.Range("C40").Value = '2.300130577309239e-05'
$ws.Range("C41").Value = '2.165383058631292e-05'
$ws.Range("C42").Value = '1.992503167789412e-05'
$ws.Range("C43").Value = '1.965957284425124e-05'
$ws.Range("C44").Value = '1.846661532542195e-05'
$ws.Range("C45").Value = '1.800917038278286e-05'
$ws.Range("C46").Value = '1.780613357076488e-05'
$ws.Range("C47").Value = '1.695633875209248e-05'
$ws.Range("C48").Value = '1.653265493906129e-05'
$ws.Range("A49").Value = 34
$ws.Range("B49").Value = '$F_{q}$'
$ws.Range("C49").Value = '1.608993773054081e-05'
$ws.Range("A50").Value = 89
$ws.Range("B50").Value = '$(\langle rr \vert rr \rangle)_{2}$'
$ws.Range("C50").Value = '1.588524886995224e-05'
$ws.Range("A51").Value = 94
$ws.Range("B51").Value = '$(\langle pq \vert rs \rangle)_{3}$'
$ws.Range("C51").Value = '1.577990099527796e-05'
$ws.Range("C52").Value = '1.487376350715456e-05'
$ws.Range("C53").Value = '1.425215336092793e-05'
$ws.Range("A54").Value = 29
$ws.Range("B54").Value = '$(\omega_{p})_{0}$'
$ws.Range("C54").Value = '1.25749973060615e-05'
$ws.Range("A55").Value = 43
$ws.Range("B55").Value = '$\eta_{s}$'
$ws.Range("C55").Value = '1.219210023417441e-05'
$ws.Range("C56").Value = '9.400822314449698e-06'
$ws.Range("C57").Value = '8.83672459988667e-06'
$ws.Range("C58").Value = '8.002184029049534e-06'
$ws.Range("C59").Value = '7.773497116770327e-06'
$ws.Range("C60").Value = '6.934287438516747e-06'
$ws.Range("C61").Value = '6.707737890902148e-06'
$ws.Range("C62").Value = '6.515665493591137e-06'
$ws.Range("C63").Value = '6.37992839657438e-06'
$ws.Range("C64").Value = '6.323176553834863e-06'
$ws.Range("C65").Value = '5.763640820931741e-06'
$ws.Range("C66").Value = '5.530928990260214e-06'
$ws.Range("C67").Value = '4.867065064817177e-06'
$ws.Range("C68").Value = '4.478877549390703e-06'
$ws.Range("A69").Value = 25
$ws.Range("B69").Value = '$type_2$'
$ws.Range("C69").Value = '4.406606867648321e-06'
$ws.Range("A70").Value = 46
$ws.Range("B70").Value = '$(F_{p})_{1}$'
$ws.Range("C70").Value = '4.356473673685131e-06'
$ws.Range("A71").Value = 64
$ws.Range("B71").Value = '$(F_{r}^{\text{SCF}})_{3}$'
$ws.Range("C71").Value = '4.109587831762429e-06'
$ws.Range("A72").Value = 59
$ws.Range("B72").Value = '$(\eta_{r})_{2}$'
$ws.Range("C72").Value = '3.719487912802046e-06'
$ws.Range("A73").Value = 96
$ws.Range("B73").Value = '$(\langle pp \vert pp \rangle)_{3}$'
$ws.Range("C73").Value = '3.67550498013907e-06'
$ws.Range("C74").Value = '3.543625495980721e-06'
$ws.Range("C75").Value = '3.416834637590532e-06'
$ws.Range("C76").Value = '3.063879636861552e-06'
$ws.Range("C77").Value = '2.508079065744602e-06'
$ws.Range("C78").Value = '2.253575967651358e-06'
$ws.Range("C79").Value = '2.114986125894094e-06'
$ws.Range("C80").Value = '1.739407909806053e-06'
$ws.Range("A81").Value = 63
$ws.Range("B81").Value = '$(\eta_{p})_{3}$'
$ws.Range("C81").Value = '1.551937754210229e-06'
$ws.Range("A82").Value = 54
$ws.Range("B82").Value = '$(F_{p})_{2}$'
$ws.Range("C82").Value = '1.51437599437738e-06'
$ws.Range("A83").Value = 55
$ws.Range("B83").Value = '$(\eta_{p})_{2}$'
$ws.Range("C83").Value = '1.505936629599036e-06'
$ws.Range("C84").Value = '7.186594433268784e-07'
$ws.Range("C85").Value = '3.209896098325893e-07'
$ws.Range("C86").Value = '7.796484509038682e-08'
$ws.Range("A87").Value = 79
$ws.Range("B87").Value = '$(\langle pq \vert sr \rangle)_{1}$'
$ws.Range("C87").Value = '6.482310261648045e-08'
$ws.Range("A88").Value = 95
$ws.Range("B88").Value = '$(\langle pq \vert sr \rangle)_{3}$'
$ws.Range("C88").Value = '4.875354434503296e-08'
$ws.Range("A89").Value = 62
$ws.Range("B89").Value = '$(F_{p})_{3}$'
$ws.Range("C89").Value = '3.127032260511813e-08'
$ws.Range("A90").Value = 87
$ws.Range("B90").Value = '$(\langle pq \vert sr \rangle)_{2}$'
$ws.Range("C90").Value = '2.650666839959123e-08'
$ws.Range("C91").Value = '2.105066786859708e-08'
$ws.Range("A92").Value = 53
$ws.Range("B92").Value = '$(\omega_{p})_{2}$'
$ws.Range("C92").Value = '1.902253489822108e-08'
$ws.Range("A93").Value = 23
$ws.Range("B93").Value = '$type_0$'
$ws.Range("C93").Value = '1.819566422397973e-08'
$ws.Range("A94").Value = 65
$ws.Range("B94").Value = '$(\omega_{r})_{3}$'
$ws.Range("C94").Value = '1.808518129591661e-08'
$ws.Range("A95").Value = 69
$ws.Range("B95").Value = '$(\langle pq \vert sr \rangle)_{0}$'
$ws.Range("C95").Value = '1.669330046559723e-08'
$ws.Range("A96").Value = 18
$ws.Range("B96").Value = '(h$_{rs}$)$_{0}$'
$ws.Range("C96").Value = '8.459537605023353e-09'
$ws.Range("A97").Value = 41
$ws.Range("B97").Value = '$\omega_{s}$'
$ws.Range("C97").Value = '7.413743962668076e-09'
$ws.Range("A98").Value = 45
$ws.Range("B98").Value = '$(\omega_{p})_{1}$'
$ws.Range("C98").Value = '4.383706751824972e-09'
$ws.Range("A99").Value = 49
$ws.Range("B99").Value = '$(\omega_{r})_{1}$'
$ws.Range("C99").Value = '3.656948125825397e-09'
$ws.Range("A100").Value = 27
$ws.Range("B100").Value = '$\mathbf{b}$'
$ws.Range("C100").Value = '3.462539316647511e-09'
$ws.Range("A101").Value = 35
$ws.Range("B101").Value = '$\eta_{q}$'
$ws.Range("C101").Value = '3.165242136945541e-09'
$ws.Range("C102").Value = '8.609154523189484e-10'
$ws.Range("C103").Value = '4.090690866067805e-10'
